$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd date in F2 ("12-03-2o11" -> "12-03-2011") and mark it as a
# quote-prefixed, date-formatted text value (matches F3's date-ish styling).
$ws.Range("F2").Value = "'12-03-2011"
$ws.Range("F2").NumberFormat = "mm-dd-yy"

# Move the active selection to F8, as last left by the author.
$ws.Range("F8").Select() | Out-Null
